# Add a "commodity" column to the "config" sheet, between "sector" and
# "level", and update the selection/cursor positions left behind on the
# "MERtoPPP" and "config" sheets, matching the author's manual edit.

$wb = $excel.ActiveWorkbook

# --- "config" sheet: insert a "commodity" column, drop the old "year" column ---
$ws = $wb.Worksheets.Item("config")

# Header row
$ws.Range("C1").Value = "commodity"
$ws.Range("D1").Value = "level"

# Data rows: commodity re-uses the same value as "sector" for this test
# fixture, and "level" keeps its original values ("useful").
$ws.Range("C2").Value = "i_therm"
$ws.Range("D2").Value = "useful"

$ws.Range("C3").Value = "i_therm"
$ws.Range("D3").Value = "useful"

# The old "year" column (2020/2030/2040 down column D, with an extra row 4)
# is gone - clear out what is now the leftover, now-unused row.
$ws.Range("D4").ClearContents()

# The user widened the new "commodity" column to fit its contents.
$ws.Columns.Item(3).ColumnWidth = 10.14

# Leave the selection where the author's cursor ended up.
$ws.Range("D4").Select()

# --- "MERtoPPP" sheet: just a leftover cursor move, no data changed ---
$ws2 = $wb.Worksheets.Item("MERtoPPP")
$ws2.Range("D13").Select()

# Re-activate "config" so it stays the active/visible sheet, as before.
$ws.Activate()
